$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Fast & Furious - Supercharged"
$ws.Range("B12").Value = 28.478000000000002
$ws.Range("C12").Value = -81.4696

$ws.Range("K13").Select()
